# Update crawl timestamp for every data row (2-514) on the single sheet,
# and update the two product rows whose rating counts changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crawl timestamp column (O) for all data rows.
$ws.Range("O2:O514").Value = "2023-01-09 15:15:25"

# Row 259 (id 3388961): ratingAmount went from 41 to 42.
$ws.Range("D259").Value = 42

# Row 458 (id 7059246): now has 1 rating worth 5 stars (was blank / 0).
$ws.Range("D458").Value = 1
$ws.Range("E458").Value = 5
